# Apply the "Add data for 2022-07-13" update to the carjacking workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet to reflect the new "through" date.
$ws.Name = "Through 2022-07-05"

# Update the header label cell (column I, row 1) that holds the
# "2022 (through 07-04)" shared string.
$ws.Range("I1").Value = "2022 (through 07-05)"

# Update the affected data cells.
$ws.Range("I8").Value = 29
$ws.Range("I14").Value = 835
